$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between C20/D20 and C21/D21
$ws.Range("C20").Value = "MEC-1NA-M.T.F."
$ws.Range("D20").Value = "-"

$ws.Range("C21").Value = "MEC-1NA-M.T.F."
$ws.Range("D21").Value = "-"
